# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the upstream data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1230
$ws1.Range("F3").Value  = 1122
$ws1.Range("F11").Value = 2223
$ws1.Range("F13").Value = 1240
$ws1.Range("F17").Value = 715
$ws1.Range("F22").Value = 4194
$ws1.Range("F25").Value = 119
$ws1.Range("F26").Value = 186
$ws1.Range("F28").Value = 602
$ws1.Range("F29").Value = 20
$ws1.Range("F32").Value = 234
$ws1.Range("F34").Value = 921
$ws1.Range("F35").Value = 118

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 776

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1230
$ws4.Range("F4").Value  = 776
$ws4.Range("F5").Value  = 1122
$ws4.Range("F16").Value = 2223
$ws4.Range("F18").Value = 1240
$ws4.Range("F23").Value = 715
$ws4.Range("F28").Value = 4194
$ws4.Range("F31").Value = 119
$ws4.Range("F32").Value = 186
$ws4.Range("F34").Value = 602
$ws4.Range("F35").Value = 20
$ws4.Range("F38").Value = 234
$ws4.Range("F40").Value = 921
$ws4.Range("F41").Value = 118
